# Insert a new "Anchorage" (ANC) colo row just above the existing
# "Adelaide" (ADL) row. This pushes ADL and every following row (through
# "Suva") down by one row, growing the used range from A1:G310 to A1:G311.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 298 (currently ADL) and everything below it down by one row.
$ws.Rows.Item(298).EntireRow.Insert()

# Populate the newly inserted row 298 with the Anchorage data.
$ws.Range("A298").Value = "ANC"
$ws.Range("B298").Value = "Anchorage, United States"
$ws.Range("C298").Value = 61.158555
$ws.Range("D298").Value = -149.890208
$ws.Range("E298").Value = "US"
$ws.Range("F298").Value = "North America"
$ws.Range("G298").Value = "Anchorage"

# EntireRow.Insert() drops the thin-border "colo code" style on column A of
# the new row; restore it by copying the formatting from the row below
# (A299, the shifted-down ADL row), which still carries the correct style.
$ws.Range("A299").Copy()
$ws.Range("A298").PasteSpecial(-4122)
